# Updates Price (D) and Volume(1h) (E) columns for the cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.651.18"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "2.370.03"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.13"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.41"
$ws.Range("E6").Value = "  -6.62%  "
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.24"
$ws.Range("E11").Value = "  -8.99%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "2.737.91"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.52"
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.92"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").Value = "2.374.23"
$ws.Range("E16").Value = "  -3.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.756"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").Value = "40.566.72"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "0.0₃0909"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  -5.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.27"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.00"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -6.39%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("E26").Value = "  -7.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.65"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -5.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.11"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.47"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "15.91"
$ws.Range("E37").Value = "  -7.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0999"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("E39").Value = "  -4.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  -8.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.82"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("D43").Value = "1.957.00"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0268"
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.56"
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.34"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("E47").Value = "  -10.11%  "
$ws.Range("D48").Value = "2.598.78"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.87"
$ws.Range("E49").Value = "  -5.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.60"
$ws.Range("E50").Value = "  -5.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.11"
$ws.Range("E51").Value = "  -4.44%  "
